$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the H1 title.
$metaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Meta description*") {
        $metaIndex = $i
        break
    }
}
if ($metaIndex -eq -1) { $metaIndex = 2 }
$metaPara = $d.Paragraphs.Item($metaIndex)
[void]$metaPara.Range.Delete()

# 2. Insert a new bold paragraph "Play Five Sound Fortune Slot for Free | Review"
#    right before the final (italic prompt) paragraph.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
[void]$lastPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($count)
$newParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Five Sound Fortune Slot for Free | Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$newPara.Range.InsertXML($newParaXml)

# 3. Replace the text of the final paragraph (keep its italic formatting) with the
#    meta-description sentence (minus the leading "Meta description: ").
$count = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($count)
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)
$finalRange.Text = "Read our review of Five Sound Fortune slot game. Play for free and discover the game's features, bonus rounds, and accessibility on different devices."

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
